# epexspot_prices.xlsx update
#  - "Prix Spot": insert a new date column ("13-nov") before column DR,
#    shifting the October columns (formerly DR:EV) one slot right (DS:EW).
#    The new column gets "-" placeholders for every hourly data row.
#  - "Gaz" / "CO2": append a new trailing row (2025-11-11) with the latest
#    day's price.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": insert date column before DR (xlShiftToRight = -4161)
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Columns("DR:DR").Insert(-4161)

$wsPrix.Range("DR1").Value = "13-nov"

for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 122).Value = "-"
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append row 150
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A150").NumberFormat = "@"
$wsGaz.Range("A150").Value = "2025-11-11"
$wsGaz.Range("A150").ClearFormats()
$wsGaz.Range("B150").Value = 28.7

# ---------------------------------------------------------------------
# Sheet "CO2": append row 150
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A150").NumberFormat = "@"
$wsCo2.Range("A150").Value = "2025-11-11"
$wsCo2.Range("A150").ClearFormats()
$wsCo2.Range("B150").Value = 80.42
